$d = $word.ActiveDocument

# Find the paragraph containing "Session timeout / logout" and remove it entirely
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Session timeout / logout*") {
        $p.Range.Delete()
        break
    }
}
